$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '70.772.96'
$ws.Cells.Item(2, 5).Value = '  +1.62%  '
$ws.Cells.Item(3, 4).Value = '3.634.30'
$ws.Cells.Item(3, 5).Value = '  +3.69%  '
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(5, 4).Value = '606.28'
$ws.Cells.Item(5, 5).Value = '  +0.26%  '
$ws.Cells.Item(6, 4).Value = '199.47'
$ws.Cells.Item(6, 5).Value = '  +2.37%  '
$ws.Cells.Item(7, 4).Value = '0.628'
$ws.Cells.Item(7, 5).Value = '  +0.41%  '
$ws.Cells.Item(8, 5).Value = '  +0.07%  '
$ws.Cells.Item(9, 5).Value = '  +9.98%  '
$ws.Cells.Item(10, 4).Value = '0.649'
$ws.Cells.Item(10, 5).Value = '  -0.05%  '
$ws.Cells.Item(11, 4).Value = '53.97'
$ws.Cells.Item(11, 5).Value = '  +1.10%  '
$ws.Cells.Item(12, 4).Value = '0.0000306'
$ws.Cells.Item(12, 5).Value = '  +1.98%  '
$ws.Cells.Item(13, 4).Value = '9.56'
$ws.Cells.Item(13, 5).Value = '  +0.70%  '
$ws.Cells.Item(14, 4).Value = '4.209.49'
$ws.Cells.Item(14, 5).Value = '  +3.59%  '
$ws.Cells.Item(15, 4).Value = '677.86'
$ws.Cells.Item(15, 5).Value = '  +14.13%  '
$ws.Cells.Item(16, 4).Value = '13.03'
$ws.Cells.Item(16, 5).Value = '  +2.52%  '
$ws.Cells.Item(17, 4).Value = '70.869.66'
$ws.Cells.Item(17, 5).Value = '  +1.51%  '
$ws.Cells.Item(18, 4).Value = '3.622.64'
$ws.Cells.Item(18, 5).Value = '  +3.44%  '
$ws.Cells.Item(19, 5).Value = '  -0.08%  '
$ws.Cells.Item(20, 5).Value = '  +0.40%  '
$ws.Cells.Item(21, 4).Value = '1.00'
$ws.Cells.Item(21, 5).Value = '  +1.30%  '
$ws.Cells.Item(22, 4).Value = '18.75'
$ws.Cells.Item(22, 5).Value = '  +3.07%  '
$ws.Cells.Item(23, 4).Value = '5.39'
$ws.Cells.Item(23, 5).Value = '  +2.22%  '
$ws.Cells.Item(24, 4).Value = '105.68'
$ws.Cells.Item(24, 5).Value = '  +4.12%  '
$ws.Cells.Item(25, 5).Value = '  -0.33%  '
$ws.Cells.Item(26, 4).Value = '3.01'
$ws.Cells.Item(26, 5).Value = '  -4.59%  '
$ws.Cells.Item(27, 4).Value = '10.46'
$ws.Cells.Item(27, 5).Value = '  -3.60%  '
$ws.Cells.Item(28, 4).Value = '9.84'
$ws.Cells.Item(28, 5).Value = '  +3.50%  '
$ws.Cells.Item(29, 4).Value = '34.14'
$ws.Cells.Item(29, 5).Value = '  +2.92%  '
$ws.Cells.Item(30, 4).Value = '4.67'
$ws.Cells.Item(30, 5).Value = '  +8.97%  '
$ws.Cells.Item(31, 4).Value = '7.17'
$ws.Cells.Item(31, 5).Value = '  +1.55%  '
$ws.Cells.Item(32, 4).Value = '12.22'
$ws.Cells.Item(32, 5).Value = '  -1.35%  '
$ws.Cells.Item(33, 4).Value = '0.116'
$ws.Cells.Item(33, 5).Value = '  +0.64%  '
$ws.Cells.Item(34, 5).Value = '  +0.46%  '
$ws.Cells.Item(35, 4).Value = '3.966.12'
$ws.Cells.Item(35, 5).Value = '  +6.56%  '
$ws.Cells.Item(36, 4).Value = '0.0₃0866'
$ws.Cells.Item(36, 5).Value = '  +5.24%  '
$ws.Cells.Item(37, 5).Value = '  -0.01%  '
$ws.Cells.Item(38, 4).Value = '3.05'
$ws.Cells.Item(38, 5).Value = '  -1.51%  '
$ws.Cells.Item(39, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(39, 4).Value = '36.86'
$ws.Cells.Item(39, 5).Value = '  +1.58%  '
$ws.Cells.Item(40, 2).Value = 'Bittensor'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(40, 4).Value = '506.12'
$ws.Cells.Item(40, 5).Value = '  +4.53%  '
$ws.Cells.Item(41, 4).Value = '0.389'
$ws.Cells.Item(41, 5).Value = '  -0.43%  '
$ws.Cells.Item(42, 4).Value = '3.55'
$ws.Cells.Item(42, 5).Value = '  -2.83%  '
$ws.Cells.Item(43, 4).Value = '0.137'
$ws.Cells.Item(43, 5).Value = '  +2.64%  '
$ws.Cells.Item(44, 2).Value = 'ThetaToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(44, 4).Value = '3.09'
$ws.Cells.Item(44, 5).Value = '  +10.00%  '
$ws.Cells.Item(45, 2).Value = 'VeChain'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(45, 4).Value = '0.0460'
$ws.Cells.Item(45, 5).Value = '  +1.76%  '
$ws.Cells.Item(47, 5).Value = '  +0.67%  '
$ws.Cells.Item(48, 4).Value = '8.69'
$ws.Cells.Item(48, 5).Value = '  +3.51%  '
$ws.Cells.Item(49, 5).Value = '  -0.32%  '
$ws.Cells.Item(50, 5).Value = '  +0.70%  '
$ws.Cells.Item(51, 2).Value = 'Mantle'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(51, 4).Value = '1.30'
$ws.Cells.Item(51, 5).Value = '  +1.55%  '
